# Dream Team Project 1
#
# Slide 10 ("Facts"), body placeholder ("Text Placeholder 2"), 2nd bullet
# paragraph currently reads (single run):
#   "Approximately 30 million participating in firearms target shoots in 2017(Statista)."
#
# The edit splits that single run into three runs and fixes "firearms" -> "firearm":
#   "Approximately 30 million participating " | "in firearm " | "target shoots in 2017(Statista)."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$para = $tr.Paragraphs(2, 1)

$full = $para.Text
$needle = "in firearms "
$idx0 = $full.IndexOf($needle)
if ($idx0 -lt 0) {
    throw "Could not locate 'in firearms ' inside paragraph: $full"
}

# TextRange.Characters(start, length) is 1-based.
$startChar = $idx0 + 1
$midRange = $para.Characters($startChar, $needle.Length)
$midRange.Text = "in firearm "

Write-Host "Final paragraph text:" $para.Text
